$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Comment-column (C) updates ---
# Row 16: note now reflects the remaining work on the "current page" style.
$ws.Range("C16").Value = "Reste : Style pour la page courante"

# Rows 18 & 20: the "Presque bon" remarks are gone - clear the notes.
$ws.Range("C18").ClearContents()
$ws.Range("C20").ClearContents()

# Rows 32 & 33: remark switched from "A faire par Vanessa" to "A FAIRE !!".
$ws.Range("C32").Value = "A FAIRE !!"
$ws.Range("C33").Value = "A FAIRE !!"

# --- Status-color (B) updates: rows 18 & 22 move from the "in progress"
#     (orange) fill to the "done" (green) fill used elsewhere in the sheet. ---
$ws.Range("B18").Interior.Color = $ws.Range("B17").Interior.Color
$ws.Range("B22").Interior.Color = $ws.Range("B17").Interior.Color

# --- View state: scroll back up and move the active selection to B16. ---
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B16").Select()
